$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos-list refresh diff.
# Column D holds numeric-looking text (e.g. "1.001", "31.034.30") that Excel
# would otherwise auto-convert to a real number, so we force text formatting,
# write the literal string, then restore the default "Normal" style so no
# stray number-format / style diff is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.034.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.959.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4872"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2951"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06956"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.985.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07805"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.498"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7007"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "281.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.046.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007771"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.210.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.539"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.517"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.863"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.192"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1050"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.641"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.570"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.472"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04922"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7562"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.172"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.734"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02011"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.706"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.560"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.135"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9029"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4455"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.119"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.013.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.433"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1256"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
